$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Day 3" test case summary numbers (rows 15-17: Written / Execution / Review)
$ws.Range("C15").Value = 2178
$ws.Range("C16").Value = 1294
$ws.Range("C17").Value = 755

# "Day 4" test case summary numbers (rows 21-23: Written / Execution / Review)
$ws.Range("C21").Value = 2360
$ws.Range("C22").Value = 1294
$ws.Range("C23").Value = 755

# Scroll the view down and move the selection to C23, matching the
# saved sheet view (topLeftCell A10 / activeCell C23).
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C23").Select()
